$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Row 3 (L002 / SSO / Dev): Execute flag flips from YES to NO
$ws.Range("A3").Value = "NO"

# Row 5 (L004 / SSO / QA): Execute flag flips from NO to YES
$ws.Range("A5").Value = "YES"

# Update the active cell selection to match the recorded cursor position
$ws.Activate()
$ws.Range("E8").Select()
